$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("RegistrationDetails")
$ws3 = $wb.Worksheets.Item("Productsdetails")

# ---- RegistrationDetails (sheet2): header row (row 1), columns B..O first ----
$ws2.Range("B1").Value = "FirstName"
$ws2.Range("C1").Value = "LastName"
$ws2.Range("D1").Value = "password"
$ws2.Range("E1").Value = "selectday"
$ws2.Range("F1").Value = "selectmont"
$ws2.Range("G1").Value = "selectyear"
$ws2.Range("H1").Value = "addressfirstname"
$ws2.Range("I1").Value = "addresslastname"
$ws2.Range("J1").Value = "address1"
$ws2.Range("K1").Value = "city"
$ws2.Range("L1").Value = "selectstate"
$ws2.Range("M1").Value = "zipcode"
$ws2.Range("N1").Value = "mobileno"
$ws2.Range("O1").Value = "addressref"

# ---- data row (row 2), columns B..O ----
$ws2.Range("B2").Value = "Sushant "
$ws2.Range("C2").Value = "Jain"
$ws2.Range("D2").Value = "password"
$ws2.Range("E2").Value = "'5"
$ws2.Range("F2").Value = "'7"
$ws2.Range("G2").Value = "'6"
$ws2.Range("H2").Value = "megha"
$ws2.Range("I2").Value = "jain"
$ws2.Range("J2").Value = "amsterdam"
$ws2.Range("K2").Value = "udaipur"
$ws2.Range("L2").Value = "'7"
$ws2.Range("M2").Value = "'20345"
$ws2.Range("N2").Value = "'0645072609"
$ws2.Range("O2").Value = "rajasthan"

# ---- column A last: email address + hyperlink ----
$ws2.Range("A1").Value = "emailid"
$ws2.Range("A2").Value = "sushant2@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:sushant2@gmail.com")
$ws2.Range("A2").Style = "Hyperlink"

# ---- column widths (best effort, bestFit-like autosizing) ----
$ws2.Columns.Item(1).ColumnWidth = 19
$ws2.Columns.Item(4).ColumnWidth = 9
$ws2.Columns.Item(5).ColumnWidth = 11
$ws2.Columns.Item(6).ColumnWidth = 10
$ws2.Columns.Item(7).ColumnWidth = 16
$ws2.Columns.Item(8).ColumnWidth = 16
$ws2.Columns.Item(13).ColumnWidth = 11
$ws2.Columns.Item(15).ColumnWidth = 10

# ---- Productsdetails (sheet3): move selection to B4, this also unsets tabSelected for sheet3 ----
$ws3.Range("B4").Select()

# ---- RegistrationDetails (sheet2): activate last so it becomes the active tab, selection C4 ----
$ws2.Range("C4").Select()
